$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "26.135.93"
    "E2" = "  -0.33%  "
    "D3" = "1.657.66"
    "E3" = "  -0.30%  "
    "E4" = "  -0.32%  "
    "D5" = "218.25"
    "E5" = "  +0.05%  "
    "D6" = "0.5278"
    "E6" = "  +0.99%  "
    "E7" = "  -0.25%  "
    "D8" = "0.2610"
    "E8" = "  -1.42%  "
    "D9" = "0.06354"
    "E9" = "  +0.98%  "
    "D10" = "20.48"
    "E10" = "  -1.90%  "
    "D11" = "0.07785"
    "E11" = "  +0.61%  "
    "D12" = "4.513"
    "D13" = "1.656.27"
    "E13" = "  -0.57%  "
    "D14" = "0.5501"
    "E14" = "  +0.92%  "
    "D15" = "0.0₅8244"
    "E15" = "  +1.06%  "
    "D16" = "65.44"
    "E16" = "  +1.29%  "
    "D17" = "26.152.76"
    "D18" = "1.003"
    "E18" = "  -0.30%  "
    "D19" = "4.586"
    "E19" = "  -1.20%  "
    "D20" = "192.43"
    "E20" = "  -0.40%  "
    "E21" = "  +0.31%  "
    "D22" = "6.038"
    "E22" = "  -0.20%  "
    "E23" = "  -0.38%  "
    "D24" = "141.93"
    "E24" = "  +1.60%  "
    "D25" = "0.1251"
    "E25" = "  +1.97%  "
    "D26" = "7.288"
    "E26" = "  +1.59%  "
    "D27" = "16.21"
    "E27" = "  +0.82%  "
    "D28" = "1.440"
    "E28" = "  +1.95%  "
    "D29" = "0.05930"
    "E29" = "  -2.02%  "
    "D30" = "1.280"
    "E30" = "  +0.13%  "
    "D31" = "3.526"
    "E31" = "  -1.62%  "
    "D32" = "3.271"
    "E32" = "  +0.25%  "
    "D33" = "1.589"
    "E33" = "  -2.13%  "
    "D34" = "0.9555"
    "E34" = "  -1.09%  "
    "D35" = "2.792"
    "E35" = "  +0.26%  "
    "D36" = "2.412"
    "E36" = "  -0.68%  "
    "D37" = "0.5716"
    "E37" = "  +1.07%  "
    "D38" = "0.01623"
    "E38" = "  +1.76%  "
    "D39" = "5.801"
    "E39" = "  -3.16%  "
    "D40" = "0.8496"
    "E40" = "  -0.72%  "
    "E41" = "  -0.20%  "
    "D42" = "103.29"
    "E42" = "  +2.98%  "
    "D43" = "1.027.58"
    "D44" = "1.802.44"
    "E44" = "  -0.11%  "
    "D45" = "57.48"
    "E45" = "  +0.72%  "
    "D46" = "1.008"
    "E46" = "  -0.26%  "
    "D47" = "1.485"
    "E47" = "  +0.52%  "
    "D49" = "0.05166"
    "E49" = "  -0.29%  "
    "D50" = "7.829"
    "E50" = "  -1.96%  "
    "D51" = "0.09734"
    "E51" = "  +0.14%  "
}

# Cells whose new value looks like a plain number (e.g. "0.5278") need the
# cell pre-formatted as Text, otherwise Excel auto-converts the typed text
# into a numeric value (losing the literal/trailing-zero formatting).
$textForceCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D16",
    "D18",
    "D19",
    "D20",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D45",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51",
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

